$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (B1, C1, D1)
$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Curve"
$ws.Range("D1").Value = "Type"

# Update the active selection to reflect the cell selected after the edit
$ws.Range("I5").Select()
